# Update SwaadSutra_Consolidated_2026-01-19.xlsx - 2026-01-19T08:37:40.360Z
# A new order (#18) came in on the "All Orders" sheet - insert it at the top
# (row 2, just under the header), pushing all existing orders down by one
# row, and roll the new order into the "Daily Summary" totals for 2026-01-19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Orders")

# Make room for the new order at the top of the data (row 2).
$ws.Rows.Item(2).Insert()

# Populate the new order row. Columns that hold digit-only / date-look-alike
# text (phone numbers, collection date, notes, etc.) are pinned to text with
# a leading apostrophe so Excel doesn't silently coerce them into numbers or
# date serials - matching how the rest of the sheet stores these columns.
$ws.Range("A2").Value = 18
$ws.Range("B2").Value = "2026-01-19 08:37"
$ws.Range("C2").Value = "Radhika Joshi"
$ws.Range("D2").Value = "C 1501"
$ws.Range("E2").Value = "'9967195227"
$ws.Range("F2").Value = "Pohe x3"
$ws.Range("G2").Value = 90
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
$ws.Range("J2").Value = "'2026-01-20"
$ws.Range("K2").Value = "08:00"
$ws.Range("L2").Value = "'"
$ws.Range("M2").Value = "'"
$ws.Range("N2").Value = "'"

# Roll the new order into the Daily Summary totals for 2026-01-19.
$summary = $wb.Worksheets.Item("Daily Summary")
$summary.Range("B2").Value = 4
$summary.Range("E2").Value = 375
$summary.Range("G2").Value = 375
